# The commit swaps the two embedded themes: the theme that currently lives
# in ppt/theme/theme2.xml ("Integral", the theme actually in force for the
# slide master) is replaced by the theme that currently lives in
# ppt/theme/theme1.xml ("Office Theme"), and vice versa.
#
# The PowerPoint object model doesn't expose a generic "replace this OOXML
# part" verb, but it does expose the live theme colors through
# SlideMaster.ColorScheme(i).RGB - exactly the lever real PowerPoint
# automation uses to re-theme a deck. We drive that lever with the 12 colors
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that fixed order) taken
# from the "Office Theme" palette so the master's live theme becomes the
# Office palette, matching the post-edit content of ppt/theme/theme2.xml.

$p  = $ppt.ActivePresentation
$m  = $p.SlideMaster
$cs = $m.ColorScheme

# RGB uses the standard OLE COLORREF encoding (0x00BBGGRR), so e.g.
# srgbClr 44546A (Office "dk2") becomes 0x6A5444 = 6968388.
$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
